$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the dated records (rows 249-251),
# pushing the existing rows 249-256 down to 252-259.
$ws.Range('A249:A251').EntireRow.Insert()

# Row 249: new weekly record - Asterix, Región Metropolitana
$ws.Range('A249').Value = 1
$ws.Range('B249').Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range('C249').Value = 'Arica y Parinacota'
$ws.Range('D249').Value = 45147
$ws.Range('E249').Value = 15
$ws.Range('F249').Value = 100114001
$ws.Range('G249').Value = 'Papa'
$ws.Range('H249').Value = 'Asterix'
$ws.Range('I249').Value = '1a (cosecha)'
$ws.Range('J249').Value = 1000
$ws.Range('K249').Value = 24000
$ws.Range('L249').Value = 25000
$ws.Range('M249').Value = 24500
$ws.Range('N249').Value = '$/saco 25 kilos'
$ws.Range('O249').Value = 'Región Metropolitana'
$ws.Range('P249').Value = 980
$ws.Range('Q249').Value = 25
$ws.Range('R249').Value = 'Hortaliza'

# Row 250: new weekly record - Asterix, Región del Maule
$ws.Range('A250').Value = 1
$ws.Range('B250').Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range('C250').Value = 'Arica y Parinacota'
$ws.Range('D250').Value = 45147
$ws.Range('E250').Value = 15
$ws.Range('F250').Value = 100114001
$ws.Range('G250').Value = 'Papa'
$ws.Range('H250').Value = 'Asterix'
$ws.Range('I250').Value = '1a (cosecha)'
$ws.Range('J250').Value = 1000
$ws.Range('K250').Value = 24000
$ws.Range('L250').Value = 25000
$ws.Range('M250').Value = 24600
$ws.Range('N250').Value = '$/saco 25 kilos'
$ws.Range('O250').Value = 'Región del Maule'
$ws.Range('P250').Value = 984
$ws.Range('Q250').Value = 25
$ws.Range('R250').Value = 'Hortaliza'

# Row 251: new weekly record - Cardinal, Región Metropolitana
$ws.Range('A251').Value = 1
$ws.Range('B251').Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range('C251').Value = 'Arica y Parinacota'
$ws.Range('D251').Value = 45147
$ws.Range('E251').Value = 15
$ws.Range('F251').Value = 100114001
$ws.Range('G251').Value = 'Papa'
$ws.Range('H251').Value = 'Cardinal'
$ws.Range('I251').Value = '1a (cosecha)'
$ws.Range('J251').Value = 1000
$ws.Range('K251').Value = 23000
$ws.Range('L251').Value = 24000
$ws.Range('M251').Value = 23600
$ws.Range('N251').Value = '$/saco 25 kilos'
$ws.Range('O251').Value = 'Región Metropolitana'
$ws.Range('P251').Value = 944
$ws.Range('Q251').Value = 25
$ws.Range('R251').Value = 'Hortaliza'
